$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped from 45206 to 45208
# for every data row (rows 2 through 224).
$range = $ws.Range("C2:C224")
$range.Value = 45208
